$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9431.556
$ws.Range("I43").Value = 6450.5
$ws.Range("J43").Value = 10283.286
$ws.Range("K43").Value = 6450.5
$ws.Range("L43").Value = 10283.286
$ws.Range("M43").Value = -6381.5
$ws.Range("N43").Value = -10421.286
$ws.Range("H74").Value = 10640.546
$ws.Range("I74").Value = 9808.4
$ws.Range("K74").Value = 9808.4
$ws.Range("M74").Value = -8872.4
$ws.Range("H77").Value = 10640.546
$ws.Range("I77").Value = 9808.4
$ws.Range("K77").Value = 49042
$ws.Range("M77").Value = -44362
$ws.Range("H92").Value = 1556.2778
$ws.Range("I92").Value = 659.6667
$ws.Range("K92").Value = 659.6667
$ws.Range("M92").Value = 588.3333
$ws.Range("H98").Value = 8757.177
$ws.Range("I98").Value = 604.9091
$ws.Range("K98").Value = 604.9091
$ws.Range("M98").Value = 893.0909
$ws.Range("H100").Value = 8525.25
$ws.Range("I100").Value = 8138.75
$ws.Range("J100").Value = 8911.75
$ws.Range("K100").Value = 8138.75
$ws.Range("L100").Value = 8911.75
$ws.Range("M100").Value = -7597.75
$ws.Range("N100").Value = -9993.75
$ws.Range("H116").Value = 15500.9
$ws.Range("I116").Value = 14667.333
$ws.Range("J116").Value = 16751.25
$ws.Range("K116").Value = 14667.333
$ws.Range("L116").Value = 16751.25
$ws.Range("M116").Value = -11225.333
$ws.Range("N116").Value = -23635.25
$ws.Range("H122").Value = 8757.177
$ws.Range("I122").Value = 604.9091
$ws.Range("K122").Value = 1814.7273
$ws.Range("M122").Value = 635.2727
$ws.Range("H132").Value = 1587.5862
$ws.Range("I132").Value = 1612.5
$ws.Range("K132").Value = 4837.5
$ws.Range("M132").Value = -2307.5
$ws.Range("H137").Value = 2268.611
$ws.Range("I137").Value = 1323.4667
$ws.Range("K137").Value = 3970.4001
$ws.Range("M137").Value = -1420.4001
$ws.Range("H138").Value = 4342.9443
$ws.Range("J138").Value = 3731.7334
$ws.Range("L138").Value = 11195.2002
$ws.Range("N138").Value = -21475.2002

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1794.94
$ws.Range("I32").Value = 1255.174
$ws.Range("K32").Value = 1255.174
$ws.Range("M32").Value = -968.174
$ws.Range("H95").Value = 39000
$ws.Range("J95").Value = 39000
$ws.Range("L95").Value = 39000
$ws.Range("N95").Value = -44492
$ws.Range("H134").Value = 116449.5
$ws.Range("J134").Value = 116449.5
$ws.Range("L134").Value = 116449.5
$ws.Range("N134").Value = -126589.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5356.5713
$ws.Range("I86").Value = 5249.3335
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 5249.3335
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -4126.3335
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5356.5713
$ws.Range("I89").Value = 5249.3335
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 26246.6675
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -20630.6675
$ws.Range("N89").Value = -41232
$ws.Range("H99").Value = 1594.9
$ws.Range("I99").Value = 1583
$ws.Range("J99").Value = 1612.75
$ws.Range("K99").Value = 1583
$ws.Range("L99").Value = 1612.75
$ws.Range("M99").Value = -85
$ws.Range("N99").Value = -4608.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3228.6428
$ws.Range("I16").Value = 2319
$ws.Range("J16").Value = 5502.75
$ws.Range("K16").Value = 2319
$ws.Range("L16").Value = 5502.75
$ws.Range("M16").Value = -2032
$ws.Range("N16").Value = -6076.75
$ws.Range("H31").Value = 33352.668
$ws.Range("I31").Value = 4408.926
$ws.Range("J31").Value = 120183.89
$ws.Range("K31").Value = 4408.926
$ws.Range("L31").Value = 120183.89
$ws.Range("M31").Value = -4113.926
$ws.Range("N31").Value = -120773.89
$ws.Range("H34").Value = 33352.668
$ws.Range("I34").Value = 4408.926
$ws.Range("J34").Value = 120183.89
$ws.Range("K34").Value = 4408.926
$ws.Range("L34").Value = 120183.89
$ws.Range("M34").Value = -4206.926
$ws.Range("N34").Value = -120587.89
$ws.Range("H99").Value = 2680.5
$ws.Range("I99").Value = 2508.5557
$ws.Range("J99").Value = 3196.3333
$ws.Range("K99").Value = 2508.5557
$ws.Range("L99").Value = 3196.3333
$ws.Range("M99").Value = -1010.5557
$ws.Range("N99").Value = -6192.3333
$ws.Range("H107").Value = 1627.4
$ws.Range("J107").Value = 1211.6
$ws.Range("L107").Value = 1211.6
$ws.Range("N107").Value = -5051.6
$ws.Range("H113").Value = 3228.6428
$ws.Range("I113").Value = 2319
$ws.Range("J113").Value = 5502.75
$ws.Range("K113").Value = 2319
$ws.Range("L113").Value = 5502.75
$ws.Range("M113").Value = -149
$ws.Range("N113").Value = -9842.75
$ws.Range("H126").Value = 2680.5
$ws.Range("I126").Value = 2508.5557
$ws.Range("J126").Value = 3196.3333
$ws.Range("K126").Value = 7525.6671
$ws.Range("L126").Value = 9588.999899999999
$ws.Range("M126").Value = -5055.6671
$ws.Range("N126").Value = -14528.9999
$ws.Range("H134").Value = 3330.1
$ws.Range("I134").Value = 2335.3635
$ws.Range("J134").Value = 4545.8887
$ws.Range("K134").Value = 7006.0905
$ws.Range("L134").Value = 13637.6661
$ws.Range("M134").Value = -4471.0905
$ws.Range("N134").Value = -18707.6661

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6243567
$ws.Range("I4").Value = 3333609.8
$ws.Range("J4").Value = 16220563
$ws.Range("K4").Value = 10000829.4
$ws.Range("L4").Value = 48661689
$ws.Range("M4").Value = -10000717.4
$ws.Range("N4").Value = -48661913
$ws.Range("H56").Value = 5667.9
$ws.Range("I56").Value = 5667.9
$ws.Range("K56").Value = 5667.9
$ws.Range("M56").Value = -5137.9
$ws.Range("H134").Value = 5491.4287
$ws.Range("I134").Value = 949.8570999999999
$ws.Range("J134").Value = 10033
$ws.Range("K134").Value = 2849.5713
$ws.Range("L134").Value = 30099
$ws.Range("M134").Value = 2220.4287
$ws.Range("N134").Value = -40239
$ws.Range("H136").Value = 2333
$ws.Range("I136").Value = 2333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1899
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 5153.4375
$ws.Range("J137").Value = 6210.8335
$ws.Range("L137").Value = 18632.5005
$ws.Range("N137").Value = -28832.5005
$ws.Range("H139").Value = 4454.722
$ws.Range("I139").Value = 1729.5454
$ws.Range("J139").Value = 8737.143
$ws.Range("K139").Value = 5188.6362
$ws.Range("L139").Value = 26211.429
$ws.Range("M139").Value = -48.63619999999992
$ws.Range("N139").Value = -36491.429

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 642.8570999999999
$ws.Range("I9").Value = 133.33333
$ws.Range("J9").Value = 1025
$ws.Range("K9").Value = 133.33333
$ws.Range("L9").Value = 1025
$ws.Range("M9").Value = 36.66667000000001
$ws.Range("N9").Value = -1365
$ws.Range("H113").Value = 5238.0527
$ws.Range("I113").Value = 4828.0835
$ws.Range("J113").Value = 5940.857
$ws.Range("K113").Value = 4828.0835
$ws.Range("L113").Value = 5940.857
$ws.Range("M113").Value = -2658.0835
$ws.Range("N113").Value = -10280.857
$ws.Range("H122").Value = 17317.092
$ws.Range("I122").Value = 21682
$ws.Range("K122").Value = 65046
$ws.Range("M122").Value = -62596
$ws.Range("H126").Value = 4322.6
$ws.Range("I126").Value = 2690.2222
$ws.Range("K126").Value = 8070.6666
$ws.Range("M126").Value = -5600.6666

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8862.933999999999
$ws.Range("I40").Value = 7994.923
$ws.Range("J40").Value = 14505
$ws.Range("K40").Value = 7994.923
$ws.Range("L40").Value = 14505
$ws.Range("M40").Value = -7858.923
$ws.Range("N40").Value = -14777
$ws.Range("H41").Value = 30000
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30876
$ws.Range("H55").Value = 3848437.2
$ws.Range("I55").Value = 8333924.5
$ws.Range("J55").Value = 3733.8572
$ws.Range("K55").Value = 8333924.5
$ws.Range("L55").Value = 3733.8572
$ws.Range("M55").Value = -8333751.5
$ws.Range("N55").Value = -4079.8572
$ws.Range("H122").Value = 9754.5
$ws.Range("I122").Value = 5004
$ws.Range("J122").Value = 14505
$ws.Range("K122").Value = 15012
$ws.Range("L122").Value = 43515
$ws.Range("M122").Value = -12562
$ws.Range("N122").Value = -48415

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6224.75
$ws.Range("J62").Value = 5666.3335
$ws.Range("L62").Value = 5666.3335
$ws.Range("N62").Value = -6914.3335
$ws.Range("H65").Value = 6224.75
$ws.Range("J65").Value = 5666.3335
$ws.Range("L65").Value = 28331.6675
$ws.Range("N65").Value = -34571.6675
$ws.Range("H70").Value = 21582.834
$ws.Range("I70").Value = 21582.834
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 21582.834
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -21267.834
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 21582.834
$ws.Range("I73").Value = 21582.834
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 21582.834
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -20490.834
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 5219.1
$ws.Range("I81").Value = 2527.4285
$ws.Range("J81").Value = 11499.667
$ws.Range("K81").Value = 5054.857
$ws.Range("L81").Value = 22999.334
$ws.Range("M81").Value = -3993.857
$ws.Range("N81").Value = -25121.334
$ws.Range("H84").Value = 5219.1
$ws.Range("I84").Value = 2527.4285
$ws.Range("J84").Value = 11499.667
$ws.Range("K84").Value = 25274.285
$ws.Range("L84").Value = 114996.67
$ws.Range("M84").Value = -19970.285
$ws.Range("N84").Value = -125604.67
$ws.Range("H126").Value = 3024.4211
$ws.Range("I126").Value = 2735.6924
$ws.Range("K126").Value = 8207.0772
$ws.Range("M126").Value = -5737.0772
$ws.Range("H136").Value = 3375.4546
$ws.Range("I136").Value = 2012.75
$ws.Range("K136").Value = 6038.25
$ws.Range("M136").Value = -3488.25
